$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the bird counts for the existing species (rows 2-8) by 1
$ws.Range("B2").Value = 15
$ws.Range("B3").Value = 14
$ws.Range("B4").Value = 13
$ws.Range("B5").Value = 12
$ws.Range("B6").Value = 10
$ws.Range("B7").Value = 9
$ws.Range("B8").Value = 9

# Add the new "Striated Heron" row at the bottom of the table
$ws.Range("A11").Value = "Striated Heron"
$ws.Range("B11").Value = 1

# Match the formatting of the species-name column used by the row above
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
